$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Cells.Item(130, 2).Value = 7483281
$ws.Cells.Item(130, 6).Value = 'SD Aucas'
$ws.Cells.Item(130, 7).Value = 'Delfin SC'
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 'D'
$ws.Cells.Item(130, 11).Value = 1.909
$ws.Cells.Item(130, 13).Value = 4.2
$ws.Cells.Item(130, 14).Value = 1.909
$ws.Cells.Item(130, 15).Value = 3.5
$ws.Cells.Item(130, 16).Value = 4
$ws.Cells.Item(130, 17).Value = -0.5
$ws.Cells.Item(130, 18).Value = 1.9
$ws.Cells.Item(130, 19).Value = 1.9
$ws.Cells.Item(130, 20).Value = 2.5
$ws.Cells.Item(130, 21).Value = 1.8
$ws.Cells.Item(130, 22).Value = 2
$ws.Cells.Item(130, 23).Value = -1
$ws.Cells.Item(130, 24).Value = 2.5
$ws.Cells.Item(130, 26).Value = -1
$ws.Cells.Item(130, 27).Value = 0.8999999999999999
$ws.Cells.Item(130, 29).Value = 1

# Row 131
$ws.Cells.Item(131, 2).Value = 7483247
$ws.Cells.Item(131, 6).Value = 'Mushuc Runa'
$ws.Cells.Item(131, 7).Value = 'Universidad Catolica del Ecuador'
$ws.Cells.Item(131, 9).Value = 2
$ws.Cells.Item(131, 10).Value = 'A'
$ws.Cells.Item(131, 11).Value = 3.25
$ws.Cells.Item(131, 12).Value = 3.2
$ws.Cells.Item(131, 13).Value = 2.25
$ws.Cells.Item(131, 14).Value = 3.5
$ws.Cells.Item(131, 15).Value = 3.25
$ws.Cells.Item(131, 16).Value = 2.1
$ws.Cells.Item(131, 17).Value = 0.5
$ws.Cells.Item(131, 18).Value = 1.775
$ws.Cells.Item(131, 19).Value = 2.025
$ws.Cells.Item(131, 21).Value = 1.9
$ws.Cells.Item(131, 22).Value = 1.9
$ws.Cells.Item(131, 24).Value = -1
$ws.Cells.Item(131, 25).Value = 1.1
$ws.Cells.Item(131, 27).Value = 1.025
$ws.Cells.Item(131, 29).Value = 0.8999999999999999

# Row 132
$ws.Cells.Item(132, 2).Value = 7483081
$ws.Cells.Item(132, 6).Value = 'Deportivo Cuenca'
$ws.Cells.Item(132, 7).Value = 'El Nacional'
$ws.Cells.Item(132, 8).Value = 1
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 'H'
$ws.Cells.Item(132, 11).Value = 2.75
$ws.Cells.Item(132, 12).Value = 3.25
$ws.Cells.Item(132, 13).Value = 2.55
$ws.Cells.Item(132, 14).Value = 3
$ws.Cells.Item(132, 15).Value = 3.3
$ws.Cells.Item(132, 16).Value = 2.3
$ws.Cells.Item(132, 17).Value = 0.25
$ws.Cells.Item(132, 18).Value = 1.825
$ws.Cells.Item(132, 19).Value = 1.975
$ws.Cells.Item(132, 20).Value = 2.75
$ws.Cells.Item(132, 21).Value = 2
$ws.Cells.Item(132, 22).Value = 1.8
$ws.Cells.Item(132, 23).Value = 2
$ws.Cells.Item(132, 25).Value = -1
$ws.Cells.Item(132, 26).Value = 0.825
$ws.Cells.Item(132, 27).Value = -1
$ws.Cells.Item(132, 29).Value = 0.8

# Row 134
$ws.Cells.Item(134, 2).Value = 7483188
$ws.Cells.Item(134, 6).Value = 'Gualaceo SC'
$ws.Cells.Item(134, 7).Value = 'Emelec'
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 2
$ws.Cells.Item(134, 10).Value = 'A'
$ws.Cells.Item(134, 11).Value = 3.6
$ws.Cells.Item(134, 12).Value = 3.3
$ws.Cells.Item(134, 13).Value = 2.05
$ws.Cells.Item(134, 14).Value = 2.6
$ws.Cells.Item(134, 15).Value = 3.25
$ws.Cells.Item(134, 16).Value = 2.75
$ws.Cells.Item(134, 17).Value = 0
$ws.Cells.Item(134, 18).Value = 1.8
$ws.Cells.Item(134, 19).Value = 2
$ws.Cells.Item(134, 20).Value = 2.5
$ws.Cells.Item(134, 21).Value = 1.975
$ws.Cells.Item(134, 22).Value = 1.825
$ws.Cells.Item(134, 24).Value = -1
$ws.Cells.Item(134, 25).Value = 1.75
$ws.Cells.Item(134, 27).Value = 1
$ws.Cells.Item(134, 28).Value = -1
$ws.Cells.Item(134, 29).Value = 0.825

# Row 135
$ws.Cells.Item(135, 2).Value = 7483306
$ws.Cells.Item(135, 6).Value = 'Tecnico Universitario'
$ws.Cells.Item(135, 7).Value = 'Club Atletico Libertad'
$ws.Cells.Item(135, 8).Value = 1
$ws.Cells.Item(135, 10).Value = 'D'
$ws.Cells.Item(135, 11).Value = 1.5
$ws.Cells.Item(135, 12).Value = 4.333
$ws.Cells.Item(135, 13).Value = 5.75
$ws.Cells.Item(135, 14).Value = 1.533
$ws.Cells.Item(135, 15).Value = 4.2
$ws.Cells.Item(135, 16).Value = 5.5
$ws.Cells.Item(135, 17).Value = -1
$ws.Cells.Item(135, 18).Value = 1.925
$ws.Cells.Item(135, 19).Value = 1.875
$ws.Cells.Item(135, 20).Value = 2.25
$ws.Cells.Item(135, 21).Value = 1.8
$ws.Cells.Item(135, 22).Value = 2
$ws.Cells.Item(135, 23).Value = -1
$ws.Cells.Item(135, 24).Value = 3.2
$ws.Cells.Item(135, 26).Value = -1
$ws.Cells.Item(135, 27).Value = 0.875
$ws.Cells.Item(135, 28).Value = -0.5
$ws.Cells.Item(135, 29).Value = 0.5

# Row 136
$ws.Cells.Item(136, 2).Value = 7482832
$ws.Cells.Item(136, 6).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(136, 7).Value = 'Guayaquil City'
$ws.Cells.Item(136, 8).Value = 2
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 10).Value = 'H'
$ws.Cells.Item(136, 11).Value = 1.363
$ws.Cells.Item(136, 12).Value = 5
$ws.Cells.Item(136, 13).Value = 7.5
$ws.Cells.Item(136, 14).Value = 1.444
$ws.Cells.Item(136, 15).Value = 4
$ws.Cells.Item(136, 16).Value = 8
$ws.Cells.Item(136, 17).Value = -1.25
$ws.Cells.Item(136, 18).Value = 2.05
$ws.Cells.Item(136, 19).Value = 1.75
$ws.Cells.Item(136, 21).Value = 1.95
$ws.Cells.Item(136, 22).Value = 1.85
$ws.Cells.Item(136, 23).Value = 0.444
$ws.Cells.Item(136, 25).Value = -1
$ws.Cells.Item(136, 26).Value = -0.5
$ws.Cells.Item(136, 27).Value = 0.375
$ws.Cells.Item(136, 28).Value = 0.95
$ws.Cells.Item(136, 29).Value = -1

# Row 139
$ws.Cells.Item(139, 2).Value = 7528859
$ws.Cells.Item(139, 6).Value = 'Club Atletico Libertad'
$ws.Cells.Item(139, 7).Value = 'Cumbaya FC'
$ws.Cells.Item(139, 8).Value = 3
$ws.Cells.Item(139, 9).Value = 1
$ws.Cells.Item(139, 10).Value = 'H'
$ws.Cells.Item(139, 11).Value = 1.727
$ws.Cells.Item(139, 13).Value = 4.333
$ws.Cells.Item(139, 14).Value = 1.4
$ws.Cells.Item(139, 15).Value = 4.2
$ws.Cells.Item(139, 16).Value = 7
$ws.Cells.Item(139, 17).Value = -1.25
$ws.Cells.Item(139, 18).Value = 2
$ws.Cells.Item(139, 19).Value = 1.8
$ws.Cells.Item(139, 21).Value = 1.95
$ws.Cells.Item(139, 22).Value = 1.85
$ws.Cells.Item(139, 23).Value = 0.3999999999999999
$ws.Cells.Item(139, 25).Value = -1
$ws.Cells.Item(139, 26).Value = 1
$ws.Cells.Item(139, 27).Value = -1
$ws.Cells.Item(139, 28).Value = 0.95
$ws.Cells.Item(139, 29).Value = -1

# Row 140
$ws.Cells.Item(140, 2).Value = 7528849
$ws.Cells.Item(140, 6).Value = 'Guayaquil City'
$ws.Cells.Item(140, 7).Value = 'Gualaceo SC'
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = 'A'
$ws.Cells.Item(140, 11).Value = 1.833
$ws.Cells.Item(140, 13).Value = 3.75
$ws.Cells.Item(140, 14).Value = 2.15
$ws.Cells.Item(140, 15).Value = 3.4
$ws.Cells.Item(140, 16).Value = 3
$ws.Cells.Item(140, 17).Value = -0.25
$ws.Cells.Item(140, 18).Value = 1.825
$ws.Cells.Item(140, 19).Value = 1.975
$ws.Cells.Item(140, 21).Value = 1.85
$ws.Cells.Item(140, 22).Value = 1.95
$ws.Cells.Item(140, 23).Value = -1
$ws.Cells.Item(140, 25).Value = 2
$ws.Cells.Item(140, 26).Value = -1
$ws.Cells.Item(140, 27).Value = 0.9750000000000001
$ws.Cells.Item(140, 28).Value = -1
$ws.Cells.Item(140, 29).Value = 0.95

# Row 142
$ws.Cells.Item(142, 2).Value = 7528858
$ws.Cells.Item(142, 6).Value = 'Orense'
$ws.Cells.Item(142, 7).Value = 'SD Aucas'
$ws.Cells.Item(142, 8).Value = 1
$ws.Cells.Item(142, 10).Value = 'A'
$ws.Cells.Item(142, 11).Value = 2.2
$ws.Cells.Item(142, 12).Value = 3.2
$ws.Cells.Item(142, 13).Value = 3.2
$ws.Cells.Item(142, 14).Value = 1.95
$ws.Cells.Item(142, 15).Value = 3.2
$ws.Cells.Item(142, 16).Value = 3.8
$ws.Cells.Item(142, 17).Value = -0.5
$ws.Cells.Item(142, 18).Value = 1.95
$ws.Cells.Item(142, 19).Value = 1.85
$ws.Cells.Item(142, 21).Value = 1.85
$ws.Cells.Item(142, 22).Value = 1.95
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = 2.8
$ws.Cells.Item(142, 26).Value = -1
$ws.Cells.Item(142, 27).Value = 0.8500000000000001
$ws.Cells.Item(142, 28).Value = 0.8500000000000001

# Row 143
$ws.Cells.Item(143, 2).Value = 7528852
$ws.Cells.Item(143, 6).Value = 'Delfin SC'
$ws.Cells.Item(143, 7).Value = 'Tecnico Universitario'
$ws.Cells.Item(143, 8).Value = 2
$ws.Cells.Item(143, 10).Value = 'D'
$ws.Cells.Item(143, 11).Value = 2.1
$ws.Cells.Item(143, 12).Value = 3.4
$ws.Cells.Item(143, 13).Value = 3.1
$ws.Cells.Item(143, 14).Value = 2.1
$ws.Cells.Item(143, 15).Value = 3.4
$ws.Cells.Item(143, 16).Value = 3.1
$ws.Cells.Item(143, 17).Value = -0.25
$ws.Cells.Item(143, 18).Value = 1.8
$ws.Cells.Item(143, 19).Value = 2
$ws.Cells.Item(143, 21).Value = 1.9
$ws.Cells.Item(143, 22).Value = 1.9
$ws.Cells.Item(143, 24).Value = 2.4
$ws.Cells.Item(143, 25).Value = -1
$ws.Cells.Item(143, 26).Value = -0.5
$ws.Cells.Item(143, 27).Value = 0.5
$ws.Cells.Item(143, 28).Value = 0.8999999999999999

# Row 144
$ws.Cells.Item(144, 2).Value = 7528857
$ws.Cells.Item(144, 6).Value = 'Universidad Catolica del Ecuador'
$ws.Cells.Item(144, 7).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(144, 8).Value = 0
$ws.Cells.Item(144, 10).Value = 'A'
$ws.Cells.Item(144, 11).Value = 1.533
$ws.Cells.Item(144, 12).Value = 4
$ws.Cells.Item(144, 13).Value = 5.5
$ws.Cells.Item(144, 14).Value = 1.5
$ws.Cells.Item(144, 15).Value = 4.333
$ws.Cells.Item(144, 16).Value = 5.25
$ws.Cells.Item(144, 17).Value = -1
$ws.Cells.Item(144, 18).Value = 1.8
$ws.Cells.Item(144, 19).Value = 2
$ws.Cells.Item(144, 20).Value = 3
$ws.Cells.Item(144, 21).Value = 1.975
$ws.Cells.Item(144, 22).Value = 1.825
$ws.Cells.Item(144, 23).Value = -1
$ws.Cells.Item(144, 25).Value = 4.25
$ws.Cells.Item(144, 26).Value = -1
$ws.Cells.Item(144, 27).Value = 1
$ws.Cells.Item(144, 28).Value = -1
$ws.Cells.Item(144, 29).Value = 0.825

# Row 145
$ws.Cells.Item(145, 2).Value = 7528848
$ws.Cells.Item(145, 6).Value = 'Emelec'
$ws.Cells.Item(145, 7).Value = 'Deportivo Cuenca'
$ws.Cells.Item(145, 8).Value = 2
$ws.Cells.Item(145, 10).Value = 'H'
$ws.Cells.Item(145, 11).Value = 1.75
$ws.Cells.Item(145, 12).Value = 3.5
$ws.Cells.Item(145, 13).Value = 4.2
$ws.Cells.Item(145, 14).Value = 2.4
$ws.Cells.Item(145, 15).Value = 3.1
$ws.Cells.Item(145, 16).Value = 2.75
$ws.Cells.Item(145, 17).Value = -0.25
$ws.Cells.Item(145, 18).Value = 2.05
$ws.Cells.Item(145, 19).Value = 1.75
$ws.Cells.Item(145, 20).Value = 2.25
$ws.Cells.Item(145, 21).Value = 1.8
$ws.Cells.Item(145, 22).Value = 2
$ws.Cells.Item(145, 23).Value = 1.4
$ws.Cells.Item(145, 25).Value = -1
$ws.Cells.Item(145, 26).Value = 1.05
$ws.Cells.Item(145, 27).Value = -1
$ws.Cells.Item(145, 28).Value = 0.8
$ws.Cells.Item(145, 29).Value = -1

# Row 192
$ws.Cells.Item(192, 14).Value = 1.833
$ws.Cells.Item(192, 18).Value = 1.825
$ws.Cells.Item(192, 19).Value = 1.975

# Row 193
$ws.Cells.Item(193, 18).Value = 1.85
$ws.Cells.Item(193, 19).Value = 1.95
$ws.Cells.Item(193, 21).Value = 1.825
$ws.Cells.Item(193, 22).Value = 1.975

# Row 194
$ws.Cells.Item(194, 14).Value = 3.75
$ws.Cells.Item(194, 16).Value = 1.95
$ws.Cells.Item(194, 18).Value = 1.825
$ws.Cells.Item(194, 19).Value = 1.975
$ws.Cells.Item(194, 21).Value = 1.925
$ws.Cells.Item(194, 22).Value = 1.875

# Row 195
$ws.Cells.Item(195, 14).Value = 1.8
$ws.Cells.Item(195, 15).Value = 3.5
$ws.Cells.Item(195, 16).Value = 4.5
$ws.Cells.Item(195, 17).Value = -0.75
$ws.Cells.Item(195, 21).Value = 1.9
$ws.Cells.Item(195, 22).Value = 1.9

# Row 196
$ws.Cells.Item(196, 18).Value = 1.975
$ws.Cells.Item(196, 19).Value = 1.825
$ws.Cells.Item(196, 21).Value = 1.95
$ws.Cells.Item(196, 22).Value = 1.85

# Row 197
$ws.Cells.Item(197, 14).Value = 2.9
$ws.Cells.Item(197, 16).Value = 2.4
$ws.Cells.Item(197, 17).Value = 0.25
$ws.Cells.Item(197, 18).Value = 1.75
$ws.Cells.Item(197, 19).Value = 2.05
$ws.Cells.Item(197, 21).Value = 1.95
$ws.Cells.Item(197, 22).Value = 1.85
